$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

Set-TextCell $ws "D2" "328.02"
Set-TextCell $ws "E2" "-1.51%"
Set-TextCell $ws "F2" "9-2-2023"
Set-TextCell $ws "G2" "0"
Set-TextCell $ws "D3" "45.14"
Set-TextCell $ws "E3" "-0.96%"
Set-TextCell $ws "F3" "9-2-2023"
Set-TextCell $ws "G3" "0"
Set-TextCell $ws "D4" "5.338"
Set-TextCell $ws "E4" "-4.60%"
Set-TextCell $ws "F4" "9-2-2023"
Set-TextCell $ws "G4" "0"
Set-TextCell $ws "D5" "0.08379"
Set-TextCell $ws "E5" "0.60%"
Set-TextCell $ws "F5" "9-2-2023"
Set-TextCell $ws "G5" "0"
Set-TextCell $ws "D6" "1.940"
Set-TextCell $ws "E6" "-4.14%"
Set-TextCell $ws "F6" "9-2-2023"
Set-TextCell $ws "G6" "0"
Set-TextCell $ws "D7" "0.9718"
Set-TextCell $ws "E7" "-1.17%"
Set-TextCell $ws "F7" "9-2-2023"
Set-TextCell $ws "G7" "0"
Set-TextCell $ws "D8" "2.507"
Set-TextCell $ws "E8" "-4.20%"
Set-TextCell $ws "F8" "9-2-2023"
Set-TextCell $ws "G8" "0"
Set-TextCell $ws "D9" "0.1113"
Set-TextCell $ws "E9" "-2.64%"
Set-TextCell $ws "F9" "9-2-2023"
Set-TextCell $ws "G9" "0"
Set-TextCell $ws "D10" "0.1917"
Set-TextCell $ws "E10" "-1.72%"
Set-TextCell $ws "F10" "9-2-2023"
Set-TextCell $ws "G10" "0"
Set-TextCell $ws "D11" "0.09646"
Set-TextCell $ws "E11" "-4.15%"
Set-TextCell $ws "F11" "9-2-2023"
Set-TextCell $ws "G11" "0"
Set-TextCell $ws "D12" "0.04630"
Set-TextCell $ws "E12" "0.83%"
Set-TextCell $ws "F12" "9-2-2023"
Set-TextCell $ws "G12" "0"
Set-TextCell $ws "D13" "0.1061"
Set-TextCell $ws "E13" "0.19%"
Set-TextCell $ws "F13" "9-2-2023"
Set-TextCell $ws "G13" "0"
Set-TextCell $ws "D14" "0.001297"
Set-TextCell $ws "E14" "2.31%"
Set-TextCell $ws "F14" "9-2-2023"
Set-TextCell $ws "G14" "0"
Set-TextCell $ws "D15" "0.006047"
Set-TextCell $ws "E15" "2.05%"
Set-TextCell $ws "F15" "9-2-2023"
Set-TextCell $ws "G15" "0"
Set-TextCell $ws "D16" "3.369"
Set-TextCell $ws "E16" "-0.09%"
Set-TextCell $ws "F16" "9-2-2023"
Set-TextCell $ws "G16" "0"
Set-TextCell $ws "D17" "4.434"
Set-TextCell $ws "E17" "0.24%"
Set-TextCell $ws "F17" "9-2-2023"
Set-TextCell $ws "G17" "0"
Set-TextCell $ws "E18" "0.72%"
Set-TextCell $ws "F18" "9-2-2023"
Set-TextCell $ws "G18" "0"
Set-TextCell $ws "D19" "8.393"
Set-TextCell $ws "E19" "-18.81%"
Set-TextCell $ws "F19" "9-2-2023"
Set-TextCell $ws "G19" "0"
Set-TextCell $ws "D20" "0.1379"
Set-TextCell $ws "E20" "0.10%"
Set-TextCell $ws "F20" "9-2-2023"
Set-TextCell $ws "G20" "0"
Set-TextCell $ws "D21" "0.2601"
Set-TextCell $ws "E21" "4.55%"
Set-TextCell $ws "F21" "9-2-2023"
Set-TextCell $ws "G21" "0"
Set-TextCell $ws "D22" "0.04177"
Set-TextCell $ws "E22" "1.46%"
Set-TextCell $ws "F22" "9-2-2023"
Set-TextCell $ws "G22" "0"
Set-TextCell $ws "D23" "0.001242"
Set-TextCell $ws "E23" "-4.48%"
Set-TextCell $ws "F23" "9-2-2023"
Set-TextCell $ws "G23" "0"
Set-TextCell $ws "D24" "0.004478"
Set-TextCell $ws "E24" "0.95%"
Set-TextCell $ws "F24" "9-2-2023"
Set-TextCell $ws "G24" "0"
Set-TextCell $ws "D25" "0.0001302"
Set-TextCell $ws "E25" "1.83%"
Set-TextCell $ws "F25" "9-2-2023"
Set-TextCell $ws "G25" "0"
Set-TextCell $ws "E26" "-20.09%"
Set-TextCell $ws "F26" "9-2-2023"
Set-TextCell $ws "G26" "0"
Set-TextCell $ws "F27" "9-2-2023"
Set-TextCell $ws "G27" "0"
Set-TextCell $ws "F28" "9-2-2023"
Set-TextCell $ws "G28" "0"
Set-TextCell $ws "F29" "9-2-2023"
Set-TextCell $ws "G29" "0"
Set-TextCell $ws "F30" "9-2-2023"
Set-TextCell $ws "G30" "0"
Set-TextCell $ws "F31" "9-2-2023"
Set-TextCell $ws "G31" "0"
Set-TextCell $ws "F32" "9-2-2023"
Set-TextCell $ws "G32" "0"
Set-TextCell $ws "F33" "9-2-2023"
Set-TextCell $ws "G33" "0"
Set-TextCell $ws "F34" "9-2-2023"
Set-TextCell $ws "G34" "0"
Set-TextCell $ws "F35" "9-2-2023"
Set-TextCell $ws "G35" "0"
Set-TextCell $ws "F36" "9-2-2023"
Set-TextCell $ws "G36" "0"
Set-TextCell $ws "F37" "9-2-2023"
Set-TextCell $ws "G37" "0"
Set-TextCell $ws "D38" "0.02752"
Set-TextCell $ws "E38" "-3.04%"
Set-TextCell $ws "F38" "9-2-2023"
Set-TextCell $ws "G38" "0"
Set-TextCell $ws "D39" "0.05647"
Set-TextCell $ws "E39" "-2.21%"
Set-TextCell $ws "F39" "9-2-2023"
Set-TextCell $ws "G39" "0"
Set-TextCell $ws "D40" "0.007803"
Set-TextCell $ws "E40" "2.01%"
Set-TextCell $ws "F40" "9-2-2023"
Set-TextCell $ws "G40" "0"
Set-TextCell $ws "D41" "0.1411"
Set-TextCell $ws "E41" "-1.78%"
Set-TextCell $ws "F41" "9-2-2023"
Set-TextCell $ws "G41" "0"
Set-TextCell $ws "D42" "0.007105"
Set-TextCell $ws "E42" "-5.07%"
Set-TextCell $ws "F42" "9-2-2023"
Set-TextCell $ws "G42" "0"
Set-TextCell $ws "D43" "0.002054"
Set-TextCell $ws "E43" "4.20%"
Set-TextCell $ws "F43" "9-2-2023"
Set-TextCell $ws "G43" "0"
Set-TextCell $ws "D44" "0.007974"
Set-TextCell $ws "E44" "-2.19%"
Set-TextCell $ws "F44" "9-2-2023"
Set-TextCell $ws "G44" "0"
Set-TextCell $ws "D45" "0.3507"
Set-TextCell $ws "F45" "9-2-2023"
Set-TextCell $ws "G45" "0"
Set-TextCell $ws "D46" "0.00006990"
Set-TextCell $ws "E46" "-2.69%"
Set-TextCell $ws "F46" "9-2-2023"
Set-TextCell $ws "G46" "0"
Set-TextCell $ws "E47" "0.32%"
Set-TextCell $ws "F47" "9-2-2023"
Set-TextCell $ws "G47" "0"
Set-TextCell $ws "D48" "0.003491"
Set-TextCell $ws "E48" "0.63%"
Set-TextCell $ws "F48" "9-2-2023"
Set-TextCell $ws "G48" "0"
Set-TextCell $ws "D49" "0.003532"
Set-TextCell $ws "E49" "1.24%"
Set-TextCell $ws "F49" "9-2-2023"
Set-TextCell $ws "G49" "0"
Set-TextCell $ws "E50" "0.32%"
Set-TextCell $ws "F50" "9-2-2023"
Set-TextCell $ws "G50" "0"
Set-TextCell $ws "E51" "0.32%"
Set-TextCell $ws "F51" "9-2-2023"
Set-TextCell $ws "G51" "0"
